$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '68.361.77'
Set-TextValue 'E2' '  +0.82%  '
Set-TextValue 'D3' '2.646.03'
Set-TextValue 'E3' '  +0.79%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '598.02'
Set-TextValue 'E5' '  +0.04%  '
Set-TextValue 'D6' '154.79'
Set-TextValue 'E6' '  +0.82%  '
Set-TextValue 'E8' '  -0.49%  '
Set-TextValue 'D9' '2.645.48'
Set-TextValue 'E9' '  +0.84%  '
Set-TextValue 'E10' '  +7.38%  '
Set-TextValue 'E11' '  -0.57%  '
Set-TextValue 'E12' '  +1.18%  '
Set-TextValue 'E13' '  +2.34%  '
Set-TextValue 'D14' '28.03'
Set-TextValue 'E14' '  +1.78%  '
Set-TextValue 'E15' '  +2.52%  '
Set-TextValue 'D16' '3.127.20'
Set-TextValue 'E16' '  +0.77%  '
Set-TextValue 'D17' '68.304.81'
Set-TextValue 'E17' '  +0.86%  '
Set-TextValue 'D18' '2.659.30'
Set-TextValue 'E18' '  +1.23%  '
Set-TextValue 'E19' '  -0.51%  '
Set-TextValue 'D20' '364.33'
Set-TextValue 'E20' '  -1.77%  '
Set-TextValue 'E22' '  +3.19%  '
Set-TextValue 'E23' '  +2.34%  '
Set-TextValue 'E24' '  +0.05%  '
Set-TextValue 'D25' '74.93'
Set-TextValue 'E25' '  +3.91%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'D27' '9.80'
Set-TextValue 'E27' '  -1.36%  '
Set-TextValue 'E30' '  +0.17%  '
Set-TextValue 'D31' '570.62'
Set-TextValue 'E31' '  -1.06%  '
Set-TextValue 'D32' '8.08'
Set-TextValue 'E32' '  +2.54%  '
Set-TextValue 'E33' '  +1.94%  '
Set-TextValue 'D34' '1.87'
Set-TextValue 'E34' '  +1.82%  '
Set-TextValue 'E35' '  +2.50%  '
Set-TextValue 'E36' '  -0.03%  '
Set-TextValue 'D37' '1.59'
Set-TextValue 'E37' '  +5.32%  '
Set-TextValue 'D38' '160.83'
Set-TextValue 'E38' '  +1.28%  '
Set-TextValue 'E39' '  +1.13%  '
Set-TextValue 'E41' '  -0.02%  '
Set-TextValue 'E42' '  +0.97%  '
Set-TextValue 'E43' '  +1.16%  '
Set-TextValue 'D44' '2.65'
Set-TextValue 'E44' '  +0.80%  '
Set-TextValue 'E45' '  +2.03%  '
Set-TextValue 'D46' '40.65'
Set-TextValue 'D48' '156.57'
Set-TextValue 'E48' '  +0.69%  '
Set-TextValue 'E49' '  +2.11%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D50' '21.97'
Set-TextValue 'E50' '  +0.38%  '
$ws.Range('B51').Value = 'Optimism'
$ws.Range('C51').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
Set-TextValue 'D51' '1.71'
Set-TextValue 'E51' '  +1.00%  '
